# Add 2022-Q3 data
#
# 1. Duplicate the existing "2022-Q2" sheet (so the new sheet inherits the
#    same layout/styling), place the copy immediately before "2022-Q2",
#    rename it to "2022-Q3" and overwrite its data with the new quarter's
#    fund holdings.
# 2. Update the "总计" (summary) sheet with the new quarter row, shifting the
#    existing quarter rows down by one.
# The other quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q2) keep their
# own data unchanged - they simply move one tab to the right because of the
# newly inserted sheet.

$wb = $excel.ActiveWorkbook

$sheetTotal = $wb.Worksheets.Item(1)
$sheetQ2 = $wb.Worksheets.Item(2)

# --- 1. Insert the new "2022-Q3" sheet before the current "2022-Q2" sheet ---
$sheetQ2.Copy($sheetQ2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Row 2 - 007216 / 浙商中华预期高股息C
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "007216"
$q3.Range("C2").NumberFormat = "@"
$q3.Range("C2").Value = "浙商中华预期高股息C"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "4.40"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "88.55"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "7.53"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.3313"
$q3.Range("H2").Value = 4

# Row 3 - 007178 / 浙商中华预期高股息A
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "007178"
$q3.Range("C3").NumberFormat = "@"
$q3.Range("C3").Value = "浙商中华预期高股息A"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "2.59"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "88.55"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "7.53"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.1950"
$q3.Range("H3").Value = 4

# --- 2. Update the "总计" sheet: add the 2022-Q3 row, shift the rest down ---
$sheetTotal.Range("B2").NumberFormat = "@"
$sheetTotal.Range("B2").Value = "2022-Q3"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.53

$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("B3").NumberFormat = "@"
$sheetTotal.Range("B3").Value = "2022-Q2"
$sheetTotal.Range("C3").Value = 2
$sheetTotal.Range("D3").Value = 1.02

$sheetTotal.Range("A4").Value = 2
$sheetTotal.Range("B4").NumberFormat = "@"
$sheetTotal.Range("B4").Value = "2022-Q1"
$sheetTotal.Range("C4").Value = 2
$sheetTotal.Range("D4").Value = 0.92

$sheetTotal.Range("A5").Value = 3
$sheetTotal.Range("B5").NumberFormat = "@"
$sheetTotal.Range("B5").Value = "2021-Q4"
$sheetTotal.Range("C5").Value = 2
$sheetTotal.Range("D5").Value = 0.23

$sheetTotal.Range("A5").Copy()
$sheetTotal.Range("A6").PasteSpecial(-4122)
$sheetTotal.Range("A6").Value = 4
$sheetTotal.Range("B6").NumberFormat = "@"
$sheetTotal.Range("B6").Value = "2021-Q2"
$sheetTotal.Range("C6").Value = 1
$sheetTotal.Range("D6").Value = 0.01

# Restore the originally active tab ("总计").
$sheetTotal.Activate()
